$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E4").Value = "2016-03-21 16:45:24"
$wsZhCn.Range("H4").Value = "2016-03-21 16:45:58"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E4").Value = "2016-03-21 16:45:29"
$wsDeDe.Range("H4").Value = "2016-03-21 16:46:18"
